$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "59.503.94"
$ws.Cells.Item(2, 5).Value = "  +0.92%  "
$ws.Cells.Item(3, 4).Value = "2.294.01"
$ws.Cells.Item(3, 5).Value = "  -0.72%  "
$ws.Cells.Item(4, 5).Value = "  +0.00%  "
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "538.94"
$ws.Cells.Item(5, 5).Value = "  -0.15%  "
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "128.53"
$ws.Cells.Item(6, 5).Value = "  -1.98%  "
$ws.Cells.Item(7, 5).Value = "  +0.02%  "
$ws.Cells.Item(8, 5).Value = "  -2.41%  "
$ws.Cells.Item(9, 4).Value = "2.289.88"
$ws.Cells.Item(9, 5).Value = "  -0.68%  "
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "0.0999"
$ws.Cells.Item(10, 5).Value = "  +0.14%  "
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "5.52"
$ws.Cells.Item(11, 5).Value = "  +0.64%  "
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "0.149"
$ws.Cells.Item(12, 5).Value = "  -0.45%  "
$ws.Cells.Item(14, 2).Value = "Avalanche"
$ws.Cells.Item(14, 3).Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "23.04"
$ws.Cells.Item(14, 5).Value = "  -2.65%  "
$ws.Cells.Item(15, 2).Value = "WrappedliquidstakedEther2.0"
$ws.Cells.Item(15, 3).Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Cells.Item(15, 4).Value = "2.703.89"
$ws.Cells.Item(15, 5).Value = "  -0.68%  "
$ws.Cells.Item(16, 4).Value = "59.478.47"
$ws.Cells.Item(16, 5).Value = "  +0.91%  "
$ws.Cells.Item(17, 5).Value = "  -0.75%  "
$ws.Cells.Item(18, 4).Value = "2.297.47"
$ws.Cells.Item(18, 5).Value = "  -0.65%  "
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "10.38"
$ws.Cells.Item(19, 5).Value = "  -1.49%  "
$ws.Cells.Item(20, 5).Value = "  -3.56%  "
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "308.29"
$ws.Cells.Item(21, 5).Value = "  -1.44%  "
$ws.Cells.Item(22, 5).Value = "  -1.66%  "
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "0.999"
$ws.Cells.Item(23, 5).Value = "  -0.15%  "
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "62.91"
$ws.Cells.Item(24, 5).Value = "  +0.53%  "
$ws.Cells.Item(25, 5).Value = "  -2.63%  "
$ws.Cells.Item(26, 5).Value = "  +0.10%  "
$ws.Cells.Item(27, 5).Value = "  -3.01%  "
$ws.Cells.Item(28, 5).Value = "  +2.57%  "
$ws.Cells.Item(29, 2).Value = "SuiNetwork"
$ws.Cells.Item(29, 3).Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "1.19"
$ws.Cells.Item(29, 5).Value = "  +2.40%  "
$ws.Cells.Item(30, 2).Value = "Monero"
$ws.Cells.Item(30, 3).Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "171.76"
$ws.Cells.Item(30, 5).Value = "  +0.31%  "
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "1.70"
$ws.Cells.Item(31, 5).Value = "  -0.90%  "
$ws.Cells.Item(32, 5).Value = "  -2.12%  "
$ws.Cells.Item(33, 5).Value = "  -1.57%  "
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "0.375"
$ws.Cells.Item(34, 5).Value = "  -2.38%  "
$ws.Cells.Item(35, 5).Value = "  +0.01%  "
$ws.Cells.Item(36, 5).Value = "  -6.67%  "
$ws.Cells.Item(37, 5).Value = "  -1.38%  "
$ws.Cells.Item(38, 5).Value = "  -0.03%  "
$ws.Cells.Item(39, 5).Value = "  -1.81%  "
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "310.70"
$ws.Cells.Item(40, 5).Value = "  -0.42%  "
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "37.37"
$ws.Cells.Item(41, 5).Value = "  -1.35%  "
$ws.Cells.Item(42, 5).Value = "  -1.03%  "
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "135.75"
$ws.Cells.Item(43, 5).Value = "  -4.07%  "
$ws.Cells.Item(44, 5).Value = "  -0.78%  "
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "0.0938"
$ws.Cells.Item(45, 5).Value = "  -1.38%  "
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "0.561"
$ws.Cells.Item(46, 5).Value = "  +0.77%  "
$ws.Cells.Item(47, 5).Value = "  +1.55%  "
$ws.Cells.Item(48, 5).Value = "  -1.35%  "
$ws.Cells.Item(49, 4).Value = "0.0₆0221"
$ws.Cells.Item(49, 5).Value = "  +21.70%  "
$ws.Cells.Item(50, 5).Value = "  +1.01%  "
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "10.99"
$ws.Cells.Item(51, 5).Value = "  -0.14%  "
